# Minor updates to "2 - About R and RStudio.pptx":
#   1. Footer/date placeholders: cached "datetimeFigureOut" field text
#      11/10/2023 -> 11/13/2023 (slide master, all slide layouts, notes master)
#   2. Slide 15 body text: shorten/reword one bullet's run text.

$p = $ppt.ActivePresentation

$oldDate = "11/10/2023"
$newDate = "11/13/2023"

function Update-DateField {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master footer/date placeholder
Update-DateField $p.SlideMaster.Shapes

# Every slide layout's footer/date placeholder
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DateField $layouts.Item($L).Shapes
}

# Notes master footer/date placeholder
Update-DateField $p.NotesMaster.Shapes

# Slide 15: "Let's head to RStudio" - reword second content bullet
# NOTE: TextRange.Text exposes plain ASCII apostrophes even though the
# underlying XML run uses the typographic RIGHT SINGLE QUOTATION MARK
# (U+2019), so search with a straight apostrophe but write back the
# curly one to match the authored text exactly.
$slide15 = $p.Slides.Item(15)
$contentShape = $slide15.Shapes.Item(2)
$tr15 = $contentShape.TextFrame.TextRange

$oldBulletAscii = "We'll put what we just learned to use in RStudio and review the basics of the program. You can watch my screen and/or run the lines yourselves. We'll assess fisheries data that you use"
$newBullet = "We" + [char]0x2019 + "ll put what we just learned to use in RStudio and review the basics of the program. You can watch my screen and write the code yourselves too. "

$fullText = $tr15.Text
$idx = $fullText.IndexOf($oldBulletAscii)
if ($idx -ge 0) {
    $sub = $tr15.Characters($idx + 1, $oldBulletAscii.Length)
    $sub.Text = $newBullet
}
